# Edit the Licor log-sheet template:
#  - shade every other sample row (the "even" numbered samples) with the
#    light "Background 2" theme fill so the log sheet is easier to read
#  - extend the print area to include the chamber-conditions block at the
#    bottom of the sheet
#  - move the selection to B4

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Apply alternating-row shading (Background 2 theme color) across A:G for
# every other data row (sample numbers 2,4,6,...,18 -> rows 3,5,7,...,19).
$shadeRows = @(3,5,7,9,11,13,15,17,19)
foreach ($r in $shadeRows) {
    $rng = $ws.Range("A" + $r + ":G" + $r)
    $rng.Interior.ThemeColor = 4
}

# Grow the print area so it covers the LICOR chamber-conditions notes too.
$ws.PageSetup.PrintArea = 'A1:G25'

# Update the current selection.
$ws.Range("B4").Select()
